$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$string16 = @'
  { "type": "Polygon",
    "coordinates": [
      [ [187671,429084], [187663,429060], [187660,429050], [187683,429031], [187695,429024], [187695,429024], [187695,429025], [187697,429031], [187697,429031], [187697,429031], [187697,429032], [187698,429035], [187699,429035], [187699,429035], [187698,429035], [187699,429036], [187698,429036], [187700,429040], [187697,429041], [187697,429042], [187693,429043], [187693,429045], [187689,429046], [187690,429048], [187686,429049], [187686,429050], [187682,429051], [187679,429052], [187680,429053], [187679,429054], [187687,429079], [187682,429081], [187671,429084] ]
      ]
   },
'@

$string17 = @'
  { "type": "Polygon",
    "coordinates": [
      [ [186752,427771], [186708,427744], [186712,427737], [186718,427726], [186797,427774], [186806,427791], [186814,427806], [186817,427811], [186820,427818], [186821,427821], [186816,427824], [186819,427831], [186831,427854], [186834,427859], [186806,427811], [186797,427816], [186799,427818], [186794,427820], [186800,427831], [186781,427837], [186798,427885], [186809,427910], [186821,427934], [186837,427981], [186854,428031], [186899,428172], [186906,428191], [186907,428196], [186904,428197], [186907,428212], [186920,428270], [186922,428281], [186919,428276], [186909,428259], [186905,428240], [186900,428220], [186894,428191], [186884,428161], [186880,428151], [186870,428117], [186859,428085], [186850,428055], [186847,428047], [186846,428043], [186841,428029], [186831,428000], [186821,427972], [186813,427948], [186807,427939], [186796,427920], [186790,427910], [186784,427896], [186778,427880], [186777,427879], [186782,427877], [186789,427874], [186788,427873], [186778,427846], [186772,427834], [186773,427833], [186789,427819], [186795,427808], [186800,427801], [186780,427788], [186780,427788], [186765,427779], [186765,427780], [186752,427771] ]
      ]
   },
'@

$string18 = @'
  { "type": "Polygon",
    "coordinates": [
      [ [187017,425818], [186960,425835], [186945,425839], [186931,425784], [186931,425784], [186931,425783], [186931,425783], [186931,425782], [186932,425782], [186932,425781], [186932,425781], [186933,425780], [186933,425780], [186933,425780], [186934,425779], [186934,425779], [186953,425775], [186969,425771], [186984,425768], [187000,425766], [187017,425765], [187016,425765], [187015,425765], [187015,425765], [187014,425766], [187014,425767], [187013,425768], [187014,425769], [187018,425793], [187021,425810], [187018,425811], [187015,425811], [187016,425813], [187017,425818] ]
      ]
   },
'@

# Row 15: reuse of the existing MultiPolygon string (same shared string, style and row height as row 14)
$ws.Range("A15").Value = 14
$ws.Range("A15").NumberFormat = "#,##0"
$ws.Range("B15").WrapText = $true
$ws.Range("B15").Value = $ws.Range("B14").Value2
$ws.Rows.Item(15).RowHeight = 144

# Row 16
$ws.Range("A16").Value = 15
$ws.Range("A16").NumberFormat = "#,##0"
$ws.Range("B16").Value = $string16
$ws.Rows.Item(16).AutoFit()

# Row 17
$ws.Range("A17").Value = 16
$ws.Range("A17").NumberFormat = "#,##0"
$ws.Range("B17").Value = $string17
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("A18").Value = 17
$ws.Range("A18").NumberFormat = "#,##0"
$ws.Range("B18").Value = $string18
$ws.Rows.Item(18).AutoFit()

# Row 19: reuse of string16
$ws.Range("A19").Value = 18
$ws.Range("A19").NumberFormat = "#,##0"
$ws.Range("B19").Value = $string16
$ws.Rows.Item(19).AutoFit()

# Update visible top-left cell and selection to match the saved view state
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("B20").Select()
